$wb = $excel.ActiveWorkbook

# This commit ("Estadisticos Matutinos 15 Oct") updates the grade
# statistics for rows 9, 12, 15 and 16 on both the "1er Parcial" and
# "3er Parcial" sheets (the "2o Parcial" sheet stays fully ungraded).
$sheetNames = @("1er Parcial", "3er Parcial")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 9: Totales=32
    $ws.Range("E9").Value = 25
    $ws.Range("F9").Value = 7
    $ws.Range("G9").Value = 78.13
    $ws.Range("H9").Value = 21.88
    $ws.Range("I9").Value = 8.2
    $ws.Range("J9").Value = 7
    $ws.Range("K9").Value = 21.88

    # Row 12: Totales=33
    $ws.Range("E12").Value = 19
    $ws.Range("F12").Value = 14
    $ws.Range("G12").Value = 57.58
    $ws.Range("H12").Value = 42.42
    $ws.Range("I12").Value = 7.4
    $ws.Range("J12").Value = 14
    $ws.Range("K12").Value = 42.42

    # Row 15: Totales=33
    $ws.Range("E15").Value = 19
    $ws.Range("F15").Value = 14
    $ws.Range("G15").Value = 57.58
    $ws.Range("H15").Value = 42.42
    $ws.Range("I15").Value = 8.3
    $ws.Range("J15").Value = 14
    $ws.Range("K15").Value = 42.42

    # Row 16: Totales=23
    $ws.Range("E16").Value = 16
    $ws.Range("F16").Value = 7
    $ws.Range("G16").Value = 69.57
    $ws.Range("H16").Value = 30.43
    $ws.Range("I16").Value = 8
    $ws.Range("J16").Value = 7
    $ws.Range("K16").Value = 30.43
}
